# Fruta / hortaliza, semanal
# Update rows 2-6 (columns D, K, L, M, N, O, P, R, S) with the new weekly data.
# This is effectively a cyclic re-shuffle of the per-row data:
#   new row2 <- old row6
#   new row3 <- old row5
#   new row4 <- old row3
#   new row5 <- old row4
#   new row6 <- old row2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44313
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1194

# Row 3
$ws.Range("D3").Value = 44305
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1361

# Row 4
$ws.Range("D4").Value = 44355
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1139

# Row 5
$ws.Range("D5").Value = 44342
$ws.Range("K5").Value = "Mankaki"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1361

# Row 6
$ws.Range("D6").Value = 44301
$ws.Range("K6").Value = "Hachiya"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1139
